# Generate Report for handback
#
# - Status goes from "Ready for handoff" to "Handed back: in sync with en-US"
#   on the Overview sheet and on each language sheet.
# - Each language sheet (zh-cn, de-de) gets two new populated columns for the
#   two real data rows: E "Latest Target File" and F "Latest Handback File",
#   mirroring the existing A (source file) / C (latest handoff file) links.
# - The "Latest Handback DateTime" column (G) gets real timestamps instead of
#   the "0001-01-01 00:00:00" placeholder, for those same two rows.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: refresh the status text shown for both tracked files.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusHandedBack
$overview.Range("C2").Value = $statusHandedBack
$overview.Range("B3").Value = $statusHandedBack
$overview.Range("C3").Value = $statusHandedBack

# ---------------------------------------------------------------------
# Helper that applies the handback updates to one language sheet.
# Parameters are positional (this host's PowerShell subset does not bind
# named/hyphenated arguments on custom functions reliably):
#   1 sheet            worksheet object (zh-cn or de-de)
#   2 mdUrl1           hyperlink target for row2 "Latest Target File" (E2)
#   3 xlfUrl1          hyperlink target for row2 "Latest Handback File" (F2)
#   4 xlfName1         display text for F2
#   5 mdUrl2           hyperlink target for row3 "Latest Target File" (E3)
#   6 xlfUrl2          hyperlink target for row3 "Latest Handback File" (F3)
#   7 xlfName2         display text for F3
#   8 handbackDatetime value written to G2/G3 "Latest Handback DateTime"
# ---------------------------------------------------------------------
function Update-LanguageSheet {
    param($sheet, $mdUrl1, $xlfUrl1, $xlfName1, $mdUrl2, $xlfUrl2, $xlfName2, $handbackDatetime)

    # Row 2
    $sheet.Range("B2").Value = $statusHandedBack
    $sheet.Hyperlinks.Add($sheet.Range("E2"), $mdUrl1, "", "", "b346fae3-79a4-44b6-8ddb-2a7b56d39f61.md") | Out-Null
    $sheet.Range("E2").Font.Underline = 2
    $sheet.Range("E2").Font.Color = 15570276
    $sheet.Hyperlinks.Add($sheet.Range("F2"), $xlfUrl1, "", "", $xlfName1) | Out-Null
    $sheet.Range("F2").Font.Underline = 2
    $sheet.Range("F2").Font.Color = 15570276
    $sheet.Range("G2").Value = $handbackDatetime

    # Row 3
    $sheet.Range("B3").Value = $statusHandedBack
    $sheet.Hyperlinks.Add($sheet.Range("E3"), $mdUrl2, "", "", "bb845ad7-5ad1-43e1-939a-6bef01431e30.md") | Out-Null
    $sheet.Range("E3").Font.Underline = 2
    $sheet.Range("E3").Font.Color = 15570276
    $sheet.Hyperlinks.Add($sheet.Range("F3"), $xlfUrl2, "", "", $xlfName2) | Out-Null
    $sheet.Range("F3").Font.Underline = 2
    $sheet.Range("F3").Font.Color = 15570276
    $sheet.Range("G3").Value = $handbackDatetime
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcnMdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/156bdc426e3298f8cb17f1ad2cf668732f09bc1c/e2e/b346fae3-79a4-44b6-8ddb-2a7b56d39f61.md"
$zhcnXlfUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bebf235f098d2fdd5c2f46f05b4abd9d528e8622/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/b346fae3-79a4-44b6-8ddb-2a7b56d39f61.45e7bdf75ea2f411a16391b6519dc4f85a59c9c1.zh-cn.xlf"
$zhcnXlfName1 = "b346fae3-79a4-44b6-8ddb-2a7b56d39f61.45e7bdf75ea2f411a16391b6519dc4f85a59c9c1.zh-cn.xlf"
$zhcnMdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/156bdc426e3298f8cb17f1ad2cf668732f09bc1c/e2e/bb845ad7-5ad1-43e1-939a-6bef01431e30.md"
$zhcnXlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bebf235f098d2fdd5c2f46f05b4abd9d528e8622/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/bb845ad7-5ad1-43e1-939a-6bef01431e30.7dbdaadff3497abdb90f15eb1d862824545e603d.zh-cn.xlf"
$zhcnXlfName2 = "bb845ad7-5ad1-43e1-939a-6bef01431e30.7dbdaadff3497abdb90f15eb1d862824545e603d.zh-cn.xlf"

Update-LanguageSheet $zhcn $zhcnMdUrl1 $zhcnXlfUrl1 $zhcnXlfName1 $zhcnMdUrl2 $zhcnXlfUrl2 $zhcnXlfName2 "2016-01-28 09:49:56"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dedeMdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/156bdc426e3298f8cb17f1ad2cf668732f09bc1c/e2e/b346fae3-79a4-44b6-8ddb-2a7b56d39f61.md"
$dedeXlfUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3f37c7109ddb3d56dec2368ab42c10d81405b0a0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/b346fae3-79a4-44b6-8ddb-2a7b56d39f61.45e7bdf75ea2f411a16391b6519dc4f85a59c9c1.de-de.xlf"
$dedeXlfName1 = "b346fae3-79a4-44b6-8ddb-2a7b56d39f61.45e7bdf75ea2f411a16391b6519dc4f85a59c9c1.de-de.xlf"
$dedeMdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/156bdc426e3298f8cb17f1ad2cf668732f09bc1c/e2e/bb845ad7-5ad1-43e1-939a-6bef01431e30.md"
$dedeXlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3f37c7109ddb3d56dec2368ab42c10d81405b0a0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/bb845ad7-5ad1-43e1-939a-6bef01431e30.7dbdaadff3497abdb90f15eb1d862824545e603d.de-de.xlf"
$dedeXlfName2 = "bb845ad7-5ad1-43e1-939a-6bef01431e30.7dbdaadff3497abdb90f15eb1d862824545e603d.de-de.xlf"

Update-LanguageSheet $dede $dedeMdUrl1 $dedeXlfUrl1 $dedeXlfName1 $dedeMdUrl2 $dedeXlfUrl2 $dedeXlfName2 "2016-01-28 09:50:18"

Write-Output "Handback report generated."
